$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.335.18"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "2.650.40"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D5").Value = "'518.07"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").Value = "'146.53"
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = "  +0.46%  "

$ws.Range("D9").Value = "2.659.87"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("E10").Value = "  -2.99%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "'0.337"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").Value = "3.116.07"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").Value = "59.330.37"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").Value = "'20.95"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "2.657.85"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "'351.69"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("D20").Value = "'4.50"
$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D21").Value = "'10.38"
$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").Value = "'6.21"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").Value = "'62.20"
$ws.Range("E24").Value = "  +1.91%  "

$ws.Range("D25").Value = "'0.417"
$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("E26").Value = "  +2.31%  "

$ws.Range("D27").Value = "'0.993"
$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").Value = "0.0₃0809"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("E31").Value = "  -3.37%  "

$ws.Range("D32").Value = "'1.58"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("D33").Value = "'18.95"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").Value = "'149.39"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").Value = "'0.952"
$ws.Range("E35").Value = "  -9.41%  "

$ws.Range("D36").Value = "'4.06"
$ws.Range("E36").Value = "  +0.96%  "

$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").Value = "'0.868"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "'36.61"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("E40").Value = "  +1.45%  "

$ws.Range("D41").Value = "'3.67"
$ws.Range("E41").Value = "  -0.97%  "

$ws.Range("D42").Value = "'279.51"
$ws.Range("E42").Value = "  -2.35%  "

$ws.Range("D43").Value = "'0.0992"
$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").Value = "'19.76"
$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").Value = "'0.603"
$ws.Range("E46").Value = "  -3.44%  "

$ws.Range("D47").Value = "2.116.90"
$ws.Range("E47").Value = "  +6.77%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0529"
$ws.Range("E48").Value = "  -2.77%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.75"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "'0.0231"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("E51").Value = "  +0.69%  "
